$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0028750212807045297
$ws.Range("D2").Value = 0.13730024614666286
$ws.Range("E2").Value = 0.14857015111023864
$ws.Range("C3").Value = 0.0028108524018678038
$ws.Range("D3").Value = 0.078299336552113644
$ws.Range("E3").Value = 0.089317705070716749
$ws.Range("C4").Value = 0.003620300416860827
$ws.Range("D4").Value = 0.25350916487195935
$ws.Range("E4").Value = 0.26770051778869913
$ws.Range("C5").Value = 0.0031164258081254185
$ws.Range("D5").Value = 0.12775536777248384
$ws.Range("E5").Value = 0.13997156524766749
$ws.Range("C6").Value = 0.003652630922463864
$ws.Range("D6").Value = 0.30599167250907716
$ws.Range("E6").Value = 0.32030975900097897
$ws.Range("C7").Value = 0.0031077802910538439
$ws.Range("D7").Value = 0.17378596994204557
$ws.Range("E7").Value = 0.18596827752209802
$ws.Range("C8").Value = 0.0032124872748273106
$ws.Range("D8").Value = 0.31048166192958337
$ws.Range("E8").Value = 0.32307441264311143
$ws.Range("C9").Value = 0.0032905055556760546
$ws.Range("D9").Value = 0.20310109856594022
$ws.Range("E9").Value = 0.21599967794380354
$ws.Range("C10").Value = 0.0031818944456713191
$ws.Range("D10").Value = 0.28732241052164859
$ws.Range("E10").Value = 0.29979523924382695
$ws.Range("C11").Value = 0.0033456651606549887
$ws.Range("D11").Value = 0.21153424121612707
$ws.Range("E11").Value = 0.22464904285261722
$ws.Range("C12").Value = 0.003423840862430388
$ws.Range("D12").Value = 0.24127090379368216
$ws.Range("E12").Value = 0.25469214745158725
$ws.Range("C13").Value = 0.0040374624366509267
$ws.Range("D13").Value = 0.20596151719115149
$ws.Range("E13").Value = 0.22178812159680303
$ws.Range("C14").Value = 0.002867933792410517
$ws.Range("D14").Value = 0.15324163287314518
$ws.Range("E14").Value = 0.16448375532253923
$ws.Range("C15").Value = 0.0041813972737935225
$ws.Range("D15").Value = 0.18868097732028813
$ws.Range("E15").Value = 0.20507179743404599
$ws.Range("C16").Value = 0.0030484869475568203
$ws.Range("D16").Value = 0.057232740408626512
$ws.Range("E16").Value = 0.069182620018995014
$ws.Range("C17").Value = 0.0039541857903619545
$ws.Range("D17").Value = 0.15041104758060247
$ws.Range("E17").Value = 0.16591121265518294
$ws.Range("C18").Value = 0.003960101831185078
$ws.Range("D18").Value = -0.019260100409157894
$ws.Range("E18").Value = -0.0037367470402389713
$ws.Range("C19").Value = 0.0052176934862527766
$ws.Range("D19").Value = 0.090645778490057366
$ws.Range("E19").Value = 0.11109881601363769
